# Insert a new weekly price-record row for "Puerro" (Vega Modelo de Temuco)
# at row 106, pushing the existing rows 106:187 down to 107:188.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("106:106").Insert()

$ws.Cells.Item(106, 1).Value  = 10
$ws.Cells.Item(106, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(106, 3).Value  = "La Araucanía"
$ws.Cells.Item(106, 4).Value  = 44634
$ws.Cells.Item(106, 5).Value  = 9
$ws.Cells.Item(106, 6).Value  = 100112005
$ws.Cells.Item(106, 7).Value  = "Puerro"
$ws.Cells.Item(106, 8).Value  = "Azul de Maquehue"
$ws.Cells.Item(106, 9).Value  = "Primera"
$ws.Cells.Item(106, 10).Value = 30
$ws.Cells.Item(106, 11).Value = 12000
$ws.Cells.Item(106, 12).Value = 12000
$ws.Cells.Item(106, 13).Value = 12000
$ws.Cells.Item(106, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(106, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(106, 16).Value = 1000
$ws.Cells.Item(106, 17).Value = 12
$ws.Cells.Item(106, 18).Value = "Hortaliza"
